{"js": "// The \"COMPETENCES TECHNIQUES\" skill lines were reordered. The text of six\n// paragraphs needs to end up in this order (paragraph formatting / pPr is\n// identical across all of them, so this is implemented as a text swap in\n// place rather than physically moving paragraphs):\n//\n//   Langages : r, python, matlab, c, c++                                      (unchanged)\n//   Visualisation : excel, tableau                                            (was position 4)\n//   MLOps : spark, vba, powerbi, Git, DVC, Flask, Docker, Github Actions,\n//           Heroku, MLflow, Streamlit                                        (was position 6)\n//   Autres : rer a auber \u00e0 nos pieds, match, stimulant, qlikview              (unchanged)\n//   ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost,\n//           OpenCV, Matplotlib, Seaborn                                      (unchanged)\n//   Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis                            (was position 2)\n\nconst newOrder = [\n  \"Langages : r, python, matlab, c, c++\",\n  \"Visualisation : excel, tableau\",\n  \"MLOps : spark, vba, powerbi, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n  \"Autres : rer a auber \u00e0 nos pieds, match, stimulant, qlikview\",\n  \"ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n  \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the six target paragraphs (by their current text) in document order,\n// keeping that order as the slots to rewrite.\nconst targetTexts = new Set([\n  \"Langages : r, python, matlab, c, c++\",\n  \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\",\n  \"Autres : rer a auber \u00e0 nos pieds, match, stimulant, qlikview\",\n  \"Visualisation : excel, tableau\",\n  \"ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n  \"MLOps : spark, vba, powerbi, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n]);\n\nconst slots = [];\nfor (const p of paragraphs.items) {\n  if (targetTexts.has(p.text)) {\n    slots.push(p);\n  }\n}\n\nif (slots.length !== newOrder.length) {\n  throw new Error(\n    `Expected ${newOrder.length} skill paragraphs, found ${slots.length}`\n  );\n}\n\nfor (let i = 0; i < slots.length; i++) {\n  slots[i].getRange().insertText(newOrder[i], \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The \"COMPETENCES TECHNIQUES\" skill lines were reordered. The text of six\n# paragraphs needs to end up in this order (paragraph formatting is identical\n# across all of them, so this is implemented as a text swap in place rather\n# than physically moving paragraphs):\n#\n#   Langages : r, python, matlab, c, c++                                      (unchanged)\n#   Visualisation : excel, tableau                                            (was position 4)\n#   MLOps : spark, vba, powerbi, Git, DVC, Flask, Docker, Github Actions,\n#           Heroku, MLflow, Streamlit                                        (was position 6)\n#   Autres : rer a auber \u00e0 nos pieds, match, stimulant, qlikview              (unchanged)\n#   ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost,\n#           OpenCV, Matplotlib, Seaborn                                      (unchanged)\n#   Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis                            (was position 2)\n\n$d = $word.ActiveDocument\n\n$newOrder = @(\n  \"Langages : r, python, matlab, c, c++\",\n  \"Visualisation : excel, tableau\",\n  \"MLOps : spark, vba, powerbi, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n  \"Autres : rer a auber \u00e0 nos pieds, match, stimulant, qlikview\",\n  \"ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n  \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\"\n)\n\n$targetTexts = @(\n  \"Langages : r, python, matlab, c, c++\",\n  \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\",\n  \"Autres : rer a auber \u00e0 nos pieds, match, stimulant, qlikview\",\n  \"Visualisation : excel, tableau\",\n  \"ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n  \"MLOps : spark, vba, powerbi, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\"\n)\n\n$slots = @()\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13)\n    if ($targetTexts -contains $t) {\n        $slots += $p\n    }\n}\n\nif ($slots.Count -ne $newOrder.Count) {\n    throw \"Expected $($newOrder.Count) skill paragraphs, found $($slots.Count)\"\n}\n\nfor ($i = 0; $i -lt $slots.Count; $i++) {\n    $slots[$i].Range.Text = $newOrder[$i]\n}\n"}
